$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: swap B1/C1 so B=Display Name, C=Projects/Employee Code
$ws.Range("B1").Value = "Display Name"
$ws.Range("C1").Value = "Projects/Employee Code"

$data = @(
    ,@("__export__.res_partner_070", "Abhijit Khanvilkar")
    ,@("__export__.res_partner_018", "Abhishek Chavan")
    ,@("__export__.res_partner_050", "Aditi Pawar")
    ,@("__export__.res_partner_011", "Ajay Dube")
    ,@("__export__.res_partner_498_3d3316b6", "Akshay Patil")
    ,@("__export__.res_partner_475_7584f238", "Amaan Shaikh")
    ,@("__export__.res_partner_058", "Amol Mestry")
    ,@("__export__.res_partner_031", "Ankita Gaikwad")
    ,@("__export__.res_partner_027", "Ankita Shirkar")
    ,@("__export__.res_partner_063", "Ankur Naik")
    ,@("__export__.res_partner_492_0c55fe59", "Ankush Nikam")
    ,@("__export__.res_partner_464_5111c6b8", "Arvind Maurya")
    ,@("__export__.res_partner_497_4db063af", "Ashish Deorukhar")
    ,@("__export__.res_partner_496_a15f9697", "Balaji Deshmukh")
    ,@("__export__.res_partner_016", "Bhavna Dengada")
    ,@("__export__.res_partner_024", "Dhanashree Khedekar")
    ,@("__export__.res_partner_055", "Divya Gire")
    ,@("__export__.res_partner_023", "Gayasuddeen Idreesi")
    ,@("__export__.res_partner_042", "Geeta Tagelly")
    ,@("__export__.res_partner_068", "Gitika Advani")
    ,@("__export__.res_partner_501_eb3c88c9", "Hansraj Bari")
    ,@("__export__.res_partner_028", "Harshali Wankhede")
    ,@("__export__.res_partner_495_5f9dbfdc", "Jaanrao Narwade")
    ,@("__export__.res_partner_469_632da1a5", "Juvil Palkar")
    ,@("__export__.res_partner_026", "Kalpana Pawar")
    ,@("__export__.res_partner_483_00098e7f", "Keyur Chotaliya")
    ,@("__export__.res_partner_008", "Kishor Dhanve")
    ,@("__export__.res_partner_060", "Kunal Patil")
    ,@("__export__.res_partner_033", "Latika Indulkar")
    ,@("__export__.res_partner_013", "Leena Pawar")
    ,@("__export__.res_partner_044", "Madhushree Sapre")
    ,@("__export__.res_partner_064", "Mahesh Gupta")
    ,@("__export__.res_partner_029", "Maheshwari Naik")
    ,@("__export__.res_partner_019", "Manjiree Jadhav")
    ,@("__export__.res_partner_069", "Manoj Harlikar")
    ,@("__export__.res_partner_041", "Mayuri Shigwan")
    ,@("__export__.res_partner_006", "Mohan Subramanian")
    ,@("__export__.res_partner_065", "Monika Sakpal")
    ,@("__export__.res_partner_057", "Namrata Kale")
    ,@("__export__.res_partner_036", "Neha Gaikwad")
    ,@("__export__.res_partner_474_8cda7fc5", "Neha Pant")
    ,@("__export__.res_partner_004", "Neha Pejle")
    ,@("__export__.res_partner_012", "Neha Talekar")
    ,@("__export__.res_partner_499_996b77af", "Neha Warkhade")
    ,@("__export__.res_partner_035", "Nikita Chile")
    ,@("__export__.res_partner_015", "Nilam Hankare")
    ,@("__export__.res_partner_049", "Nimisha Kalingan")
    ,@("__export__.res_partner_043", "Nitin Avere")
    ,@("__export__.res_partner_017", "Omkar Parab")
    ,@("__export__.res_partner_476_00cf7caa", "Pinku Adhikari")
    ,@("__export__.res_partner_021", "Prajakta Mohite")
    ,@("__export__.res_partner_066", "Prakash Dhapola")
    ,@("__export__.res_partner_009", "Pranay Zadekar")
    ,@("__export__.res_partner_038", "Pushpa Patkare")
    ,@("__export__.res_partner_047", "Rajendra Navlu")
    ,@("__export__.res_partner_477_068836f4", "Rajiv Meshram")
    ,@("__export__.res_partner_001", "Rakesh Gawade")
    ,@("__export__.res_partner_484_b938f851", "ram kishore")
    ,@("__export__.res_partner_491_c52858e2", "Ranjeev Mahindru")
    ,@("__export__.res_partner_505_cb8d2430", "Ranveer Kapoor")
    ,@("__export__.res_partner_061", "Renuka Kaushik")
    ,@("__export__.res_partner_003", "Riddhi Pujari")
    ,@("__export__.res_partner_030", "Rohan Pal")
    ,@("__export__.res_partner_074", "Rupesh Chaudhari")
    ,@("__export__.res_partner_046", "Sahil Sawant")
    ,@("__export__.res_partner_053", "Sailee Dabholkar")
    ,@("__export__.res_partner_052", "Sakshi More")
    ,@("__export__.res_partner_037", "Sanghamitra Onkar")
    ,@("__export__.res_partner_071", "Sanjay Chakraborty")
    ,@("__export__.res_partner_480_7dc8c491", "Sharim Wandrick")
    ,@("__export__.res_partner_014", "Shital Chandarkar")
    ,@("__export__.res_partner_478_d16a6bf9", "Shubham Andotra")
    ,@("__export__.res_partner_007", "Shweta Mishra")
    ,@("__export__.res_partner_010", "Sohail Kazi")
    ,@("__export__.res_partner_002", "Suhas Jadhav")
    ,@("__export__.res_partner_056", "Sunita Dhuri")
    ,@("__export__.res_partner_022", "Surekha Hande")
    ,@("__export__.res_partner_020", "Swati Gaikwad")
    ,@("__export__.res_partner_054", "Tanvi Kudalkar")
    ,@("__export__.res_partner_481_2ef310c2", "Tejas Sontakke")
    ,@("__export__.res_partner_039", "Tushar Shinde")
    ,@("__export__.res_partner_025", "Vaishnavi Khandekar")
    ,@("__export__.res_partner_040", "Varsha Shete")
    ,@("__export__.res_partner_067", "Varun Bhatia")
    ,@("__export__.res_partner_051", "Vedant Vaidya")
    ,@("__export__.res_partner_005", "Vijay Mahajan")
    ,@("__export__.res_partner_045", "Vinayak Gholap")
    ,@("__export__.res_partner_048", "Vishakha Sable")
    ,@("__export__.res_partner_062", "Vishkha Chamare")
    ,@("__export__.res_partner_032", "Vishvanath Bhosale")
    ,@("__export__.res_partner_059", "Yash Prabhulkar")
    ,@("__export__.res_partner_034", "Yogesh Kharat")
    ,@("__export__.res_partner_482_59c2d05d", "Zaid Borkar")
    ,@("__export__.res_partner_479_343c80c3", "Zeel Shah")
)

$row = 2
foreach ($pair in $data) {
    $extId = $pair[0]
    $name = $pair[1]
    $ws.Cells.Item($row, 1).Value = $extId
    $ws.Cells.Item($row, 2).Value = $name
    if ($name -eq "Abhijit Khanvilkar") {
        $ws.Cells.Item($row, 3).Value = "DRP/2024/PU0001"
    } else {
        $ws.Cells.Item($row, 3).Value = $null
    }
    $row = $row + 1
}

$lastRow = $row - 1
$ws.Range("A2:C" + $lastRow).WrapText = $true